$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows a brand-new data row inserted right before the old row 700,
# pushing the former rows 700:813 down to 701:814 (dimension grows from
# A1:R813 to A1:R814). Insert a whole row at position 700 so everything
# below shifts down and the date-format style of column D carries through.
$ws.Rows.Item(700).Insert()

# Populate the newly inserted row 700 with its data values.
$ws.Cells.Item(700, 1).Value = 6
$ws.Cells.Item(700, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(700, 3).Value = "Metropolitana"
$ws.Cells.Item(700, 4).Value = 45180
$ws.Cells.Item(700, 5).Value = 13
$ws.Cells.Item(700, 6).Value = 100112044
$ws.Cells.Item(700, 7).Value = "Perejil"
$ws.Cells.Item(700, 8).Value = "Sin especificar"
$ws.Cells.Item(700, 9).Value = "Primera"
$ws.Cells.Item(700, 10).Value = 240
$ws.Cells.Item(700, 11).Value = 11000
$ws.Cells.Item(700, 12).Value = 12000
$ws.Cells.Item(700, 13).Value = 11458
$ws.Cells.Item(700, 14).Value = "$/docena de atados"
$ws.Cells.Item(700, 15).Value = "Región Metropolitana"
$ws.Cells.Item(700, 16).Value = 3819
$ws.Cells.Item(700, 17).Value = 3
$ws.Cells.Item(700, 18).Value = "Hortaliza"
